$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('I3').Value = '%'
$ws.Range('J3').Value = 'Uninterpretable'
$ws.Range('I9').Value = 'sd'
$ws.Range('J9').Value = 'Statement-non-opinion'
$ws.Range('I14').Value = 'sv'
$ws.Range('J14').Value = 'Statement-opinion'
$ws.Range('I23').Value = 'aa'
$ws.Range('J23').Value = 'Agree/Accept'
$ws.Range('I24').Value = 'sd'
$ws.Range('J24').Value = 'Statement-non-opinion'
$ws.Range('I32').Value = 'sd'
$ws.Range('J32').Value = 'Statement-non-opinion'
$ws.Range('I33').Value = 'aa'
$ws.Range('J33').Value = 'Agree/Accept'
$ws.Range('I40').Value = 'aa'
$ws.Range('J40').Value = 'Agree/Accept'
$ws.Range('I49').Value = 'sv'
$ws.Range('J49').Value = 'Statement-opinion'
$ws.Range('I60').Value = 'aa'
$ws.Range('J60').Value = 'Agree/Accept'
$ws.Range('I62').Value = '%'
$ws.Range('J62').Value = 'Uninterpretable'
$ws.Range('I76').Value = 'sd'
$ws.Range('J76').Value = 'Statement-non-opinion'
$ws.Range('I79').Value = 'aa'
$ws.Range('J79').Value = 'Agree/Accept'
$ws.Range('I93').Value = 'sd'
$ws.Range('J93').Value = 'Statement-non-opinion'
$ws.Range('I98').Value = 'aa'
$ws.Range('J98').Value = 'Agree/Accept'
$ws.Range('I99').Value = '%'
$ws.Range('J99').Value = 'Uninterpretable'
$ws.Range('I104').Value = '%'
$ws.Range('J104').Value = 'Uninterpretable'
$ws.Range('I123').Value = 'sv'
$ws.Range('J123').Value = 'Statement-opinion'
$ws.Range('I128').Value = 'aa'
$ws.Range('J128').Value = 'Agree/Accept'
$ws.Range('I141').Value = 'aa'
$ws.Range('J141').Value = 'Agree/Accept'
$ws.Range('I154').Value = 'ba'
$ws.Range('J154').Value = 'Appreciation'
$ws.Range('I169').Value = 'sd'
$ws.Range('J169').Value = 'Statement-non-opinion'
$ws.Range('I174').Value = 'aa'
$ws.Range('J174').Value = 'Agree/Accept'
$ws.Range('I175').Value = 'sv'
$ws.Range('J175').Value = 'Statement-opinion'
$ws.Range('I181').Value = 'sd'
$ws.Range('J181').Value = 'Statement-non-opinion'
$ws.Range('I186').Value = 'sd'
$ws.Range('J186').Value = 'Statement-non-opinion'
$ws.Range('I187').Value = 'sd'
$ws.Range('J187').Value = 'Statement-non-opinion'
$ws.Range('I204').Value = '%'
$ws.Range('J204').Value = 'Uninterpretable'
$ws.Range('I222').Value = 'sd'
$ws.Range('J222').Value = 'Statement-non-opinion'
$ws.Range('I223').Value = '%'
$ws.Range('J223').Value = 'Uninterpretable'
$ws.Range('I228').Value = 'sv'
$ws.Range('J228').Value = 'Statement-opinion'
$ws.Range('I234').Value = 'aa'
$ws.Range('J234').Value = 'Agree/Accept'
$ws.Range('I235').Value = 'aa'
$ws.Range('J235').Value = 'Agree/Accept'
$ws.Range('I237').Value = 'aa'
$ws.Range('J237').Value = 'Agree/Accept'
$ws.Range('I246').Value = 'sd'
$ws.Range('J246').Value = 'Statement-non-opinion'
